$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, shifting existing rows 53..152 down to 54..153
$ws.Rows.Item(53).Insert()

# New row 53 keeps same static columns as its neighbors (A,B,C,E,F,G,H,I,N,Q,R)
# and gets new values for D, J, K, L, M, O, P.
$ws.Range("A53").Value = 11
$ws.Range("B53").Value = "Vega Monumental Concepción"
$ws.Range("C53").Value = "Bíobío"
$ws.Range("D53").NumberFormat = $ws.Range("D54").NumberFormat
$ws.Range("D53").Value = 45044
$ws.Range("E53").Value = 8
$ws.Range("F53").Value = 100112001
$ws.Range("G53").Value = "Berenjena"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 180
$ws.Range("K53").Value = 8000
$ws.Range("L53").Value = 9000
$ws.Range("M53").Value = 8556
$ws.Range("N53").Value = "$/caja 60 unidades"
$ws.Range("O53").Value = "Región Metropolitana"
$ws.Range("P53").Value = 143
$ws.Range("Q53").Value = 60
$ws.Range("R53").Value = "Hortaliza"
